# Updates the cryptos list data (prices & 1h volume %) and restores the
# original relative ordering for two coin pairs that had been re-sorted
# (Filecoin / InternetComputer(DFINITY) and ARBITRUM / NEARProtocol),
# as captured by the upstream data-refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.722.31"
$ws.Range("E2").Value = "  +4.32%  "
# Row 3
$ws.Range("D3").Value = "2.271.36"
$ws.Range("E3").Value = "  +1.82%  "
# Row 4
$ws.Range("E4").Value = "  +0.00%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.32"
$ws.Range("E5").Value = "  -0.13%  "
# Row 6
$ws.Range("E6").Value = "  +0.52%  "
# Row 7
$ws.Range("E7").Value = "  +0.74%  "
# Row 8
$ws.Range("E8").Value = "  +0.04%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.424"
$ws.Range("E9").Value = "  +5.73%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0942"
$ws.Range("E10").Value = "  +5.78%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.73"
$ws.Range("E11").Value = "  -2.04%  "
# Row 12
$ws.Range("E12").Value = "  +0.66%  "
# Row 13
$ws.Range("D13").Value = "2.611.95"
$ws.Range("E13").Value = "  +2.00%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.78"
$ws.Range("E14").Value = "  +0.95%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.75"
$ws.Range("E15").Value = "  +9.24%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.82"
$ws.Range("E16").Value = "  +4.31%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.811"
$ws.Range("E17").Value = "  +1.57%  "
# Row 18
$ws.Range("D18").Value = "2.282.17"
$ws.Range("E18").Value = "  +2.02%  "
# Row 19
$ws.Range("D19").Value = "43.707.75"
$ws.Range("E19").Value = "  +4.61%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0938"
$ws.Range("E20").Value = "  +5.33%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.12"
$ws.Range("E21").Value = "  +1.07%  "
# Row 22
$ws.Range("E22").Value = "  +3.56%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.83"
$ws.Range("E23").Value = "  +0.69%  "
# Row 24
$ws.Range("E24").Value = "  -0.07%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("E25").Value = "  +7.43%  "
# Row 26
$ws.Range("E26").Value = "  +2.27%  "
# Row 27
$ws.Range("E27").Value = "  +1.75%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.20"
$ws.Range("E28").Value = "  +2.51%  "
# Row 29
$ws.Range("E29").Value = "  -1.71%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.59"
$ws.Range("E30").Value = "  +3.23%  "
# Row 31
$ws.Range("E31").Value = "  +4.14%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.67"
$ws.Range("E32").Value = "  +1.19%  "
# Row 33
$ws.Range("E33").Value = "  +0.23%  "
# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.07"
$ws.Range("E34").Value = "  +2.43%  "
# Row 35
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.78"
$ws.Range("E35").Value = "  +3.39%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0663"
$ws.Range("E36").Value = "  +5.13%  "
# Row 37
$ws.Range("E37").Value = "  -2.23%  "
# Row 38
$ws.Range("E38").Value = "  +2.43%  "
# Row 39
$ws.Range("E39").Value = "  -1.82%  "
# Row 40
$ws.Range("E40").Value = "  +4.10%  "
# Row 41
$ws.Range("E41").Value = "  +0.07%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.78"
$ws.Range("E42").Value = "  +2.49%  "
# Row 43
$ws.Range("E43").Value = "  -11.88%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0993"
$ws.Range("E44").Value = "  +1.46%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.50"
$ws.Range("E45").Value = "  -6.53%  "
# Row 46
$ws.Range("E46").Value = "  +0.23%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.00"
$ws.Range("E47").Value = "  -0.92%  "
# Row 48
$ws.Range("D48").Value = "1.473.50"
$ws.Range("E48").Value = "  +0.03%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.63"
$ws.Range("E49").Value = "  +0.70%  "
# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.28"
$ws.Range("E50").Value = "  +8.92%  "
# Row 51
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.08"
$ws.Range("E51").Value = "  +0.67%  "
